$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing H/I values for rows with revised AgTests/AgPosit figures ---
$ws.Range("H286").Value = 54338

$ws.Range("H287").Value = 57852

$ws.Range("H288").Value = 56628
$ws.Range("I288").Value = 3989

$ws.Range("H289").Value = 65141
$ws.Range("I289").Value = 3755

$ws.Range("H292").Value = 82191
$ws.Range("I292").Value = 7298

$ws.Range("H293").Value = 82934
$ws.Range("I293").Value = 5859

$ws.Range("H294").Value = 92249
$ws.Range("I294").Value = 5106

$ws.Range("H299").Value = 65331
$ws.Range("I299").Value = 6841

$ws.Range("H300").Value = 70987
$ws.Range("I300").Value = 6949

$ws.Range("H301").Value = 70116
$ws.Range("I301").Value = 5571

$ws.Range("H302").Value = 72962
$ws.Range("I302").Value = 5329

$ws.Range("H306").Value = 70279
$ws.Range("I306").Value = 7150

$ws.Range("H307").Value = 73273
$ws.Range("I307").Value = 6338

$ws.Range("H308").Value = 17044
$ws.Range("I308").Value = 1425

$ws.Range("H309").Value = 57500
$ws.Range("I309").Value = 3996

$ws.Range("H310").Value = 90010
$ws.Range("I310").Value = 5382

$ws.Range("H313").Value = 72877
$ws.Range("I313").Value = 3546

$ws.Range("H314").Value = 64970
$ws.Range("I314").Value = 3344

$ws.Range("H315").Value = 65840
$ws.Range("I315").Value = 3095

$ws.Range("H316").Value = 49427
$ws.Range("I316").Value = 2300

# --- Append new daily rows 317-319 ---
$ws.Range("A317").Value = 44211
$ws.Range("A317").NumberFormat = "yyyy-mm-dd"
$ws.Range("B317").Value = 222752
$ws.Range("C317").Value = 166555
$ws.Range("D317").Value = 52780
$ws.Range("E317").Value = 11392
$ws.Range("F317").Value = 2045
$ws.Range("G317").Value = 3417
$ws.Range("H317").Value = 58918
$ws.Range("I317").Value = 2079

$ws.Range("A318").Value = 44212
$ws.Range("A318").NumberFormat = "yyyy-mm-dd"
$ws.Range("B318").Value = 223325
$ws.Range("C318").Value = 168915
$ws.Range("D318").Value = 50936
$ws.Range("E318").Value = 2850
$ws.Range("F318").Value = 573
$ws.Range("G318").Value = 3474
$ws.Range("H318").Value = 19705
$ws.Range("I318").Value = 783

$ws.Range("A319").Value = 44213
$ws.Range("A319").NumberFormat = "yyyy-mm-dd"
$ws.Range("B319").Value = 224385
$ws.Range("C319").Value = 171092
$ws.Range("D319").Value = 49767
$ws.Range("E319").Value = 6252
$ws.Range("F319").Value = 1060
$ws.Range("G319").Value = 3526
$ws.Range("H319").Value = 37771
$ws.Range("I319").Value = 676
